# Apply updated Spain provincias COVID data (case counts refreshed to 18:16)
# and re-sort rows by "Casos totales" (column B) descending, matching the
# source data feed ordering after the update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 18:16"

# Full data table (Ciudad, Casos totales, Casos activos, Recuperados, Muertes)
# rewritten row by row in the resulting sorted order.
$ws.Range("A4").Value = "Madrid"
$ws.Range("B4").Value = 9702
$ws.Range("C4").Value = 1186
$ws.Range("D4").Value = 6931
$ws.Range("E4").Value = 1201
$ws.Range("A5").Value = "Cataluña"
$ws.Range("B5").Value = 4704
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4078
$ws.Range("E5").Value = 191
$ws.Range("A6").Value = "Araba/Alava"
$ws.Range("B6").Value = 912
$ws.Range("C6").Value = 21
$ws.Range("D6").Value = 852
$ws.Range("E6").Value = 60
$ws.Range("A7").Value = "Bizkaia/Vizcaya"
$ws.Range("B7").Value = 860
$ws.Range("C7").Value = 21
$ws.Range("D7").Value = 836
$ws.Range("E7").Value = 24
$ws.Range("A8").Value = "Navarra"
$ws.Range("B8").Value = 794
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 652
$ws.Range("E8").Value = 14
$ws.Range("A9").Value = "La Rioja"
$ws.Range("B9").Value = 654
$ws.Range("C9").Value = 13
$ws.Range("D9").Value = 536
$ws.Range("E9").Value = 18
$ws.Range("A10").Value = "Valencia/Valencia"
$ws.Range("B10").Value = 627
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 590
$ws.Range("E10").Value = 25
$ws.Range("A11").Value = "Asturias"
$ws.Range("B11").Value = 545
$ws.Range("C11").Value = 12
$ws.Range("D11").Value = 467
$ws.Range("E11").Value = 10
$ws.Range("A12").Value = "Malaga"
$ws.Range("B12").Value = 505
$ws.Range("C12").Value = 72
$ws.Range("D12").Value = 464
$ws.Range("E12").Value = 21
$ws.Range("A13").Value = "Ciudad Real"
$ws.Range("B13").Value = 505
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 457
$ws.Range("E13").Value = 40
$ws.Range("A14").Value = "Toledo"
$ws.Range("B14").Value = 501
$ws.Range("C14").Value = 22
$ws.Range("D14").Value = 451
$ws.Range("E14").Value = 28
$ws.Range("A15").Value = "Albacete"
$ws.Range("B15").Value = 430
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 390
$ws.Range("E15").Value = 32
$ws.Range("A16").Value = "A Coruña"
$ws.Range("B16").Value = 384
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 326
$ws.Range("E16").Value = 5
$ws.Range("A17").Value = "Alacant/Alicante"
$ws.Range("B17").Value = 372
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 325
$ws.Range("E17").Value = 40
$ws.Range("A18").Value = "Pontevedra"
$ws.Range("B18").Value = 348
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 262
$ws.Range("E18").Value = 2
$ws.Range("A19").Value = "Granada"
$ws.Range("B19").Value = 335
$ws.Range("C19").Value = 72
$ws.Range("D19").Value = 276
$ws.Range("E19").Value = 13
$ws.Range("A20").Value = "Zaragoza"
$ws.Range("B20").Value = 329
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 315
$ws.Range("E20").Value = 14
$ws.Range("A21").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B21").Value = 325
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 312
$ws.Range("E21").Value = 13
$ws.Range("A22").Value = "Murcia"
$ws.Range("B22").Value = 296
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 213
$ws.Range("E22").Value = 1
$ws.Range("A23").Value = "Sevilla"
$ws.Range("B23").Value = 295
$ws.Range("C23").Value = 72
$ws.Range("D23").Value = 243
$ws.Range("E23").Value = 4
$ws.Range("A24").Value = "Cantabria"
$ws.Range("B24").Value = 282
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 200
$ws.Range("E24").Value = 5
$ws.Range("A25").Value = "Burgos"
$ws.Range("B25").Value = 269
$ws.Range("C25").Value = 27
$ws.Range("D25").Value = 175
$ws.Range("E25").Value = 16
$ws.Range("A26").Value = "Salamanca"
$ws.Range("B26").Value = 265
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 180
$ws.Range("E26").Value = 21
$ws.Range("A27").Value = "Guadalajara"
$ws.Range("B27").Value = 263
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 257
$ws.Range("E27").Value = 4
$ws.Range("A28").Value = "Tenerife"
$ws.Range("B28").Value = 262
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 211
$ws.Range("E28").Value = 5
$ws.Range("A29").Value = "Caceres"
$ws.Range("B29").Value = 243
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 231
$ws.Range("E29").Value = 12
$ws.Range("A30").Value = "Valladolid"
$ws.Range("B30").Value = 241
$ws.Range("C30").Value = 13
$ws.Range("D30").Value = 193
$ws.Range("E30").Value = 11
$ws.Range("A31").Value = "Mallorca"
$ws.Range("B31").Value = 210
$ws.Range("C31").Value = 10
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 8
$ws.Range("A32").Value = "Leon"
$ws.Range("B32").Value = 201
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 156
$ws.Range("E32").Value = 14
$ws.Range("A33").Value = "Cordoba"
$ws.Range("B33").Value = 180
$ws.Range("C33").Value = 72
$ws.Range("D33").Value = 140
$ws.Range("E33").Value = 3
$ws.Range("A34").Value = "Aragon"
$ws.Range("B34").Value = 174
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 163
$ws.Range("E34").Value = 11
$ws.Range("A35").Value = "Jaen"
$ws.Range("B35").Value = 171
$ws.Range("C35").Value = 72
$ws.Range("D35").Value = 136
$ws.Range("E35").Value = 3
$ws.Range("A36").Value = "Segovia"
$ws.Range("B36").Value = 157
$ws.Range("C36").Value = 9
$ws.Range("D36").Value = 140
$ws.Range("E36").Value = 20
$ws.Range("A37").Value = "Cadiz"
$ws.Range("B37").Value = 134
$ws.Range("C37").Value = 72
$ws.Range("D37").Value = 126
$ws.Range("E37").Value = 2
$ws.Range("A38").Value = "Cuenca"
$ws.Range("B38").Value = 120
$ws.Range("C38").Value = 8
$ws.Range("D38").Value = 104
$ws.Range("E38").Value = 8
$ws.Range("A39").Value = "Gran Canaria"
$ws.Range("B39").Value = 119
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 96
$ws.Range("E39").Value = 2
$ws.Range("A40").Value = "Soria"
$ws.Range("B40").Value = 119
$ws.Range("C40").Value = 5
$ws.Range("D40").Value = 71
$ws.Range("E40").Value = 8
$ws.Range("A41").Value = "Avila"
$ws.Range("B41").Value = 114
$ws.Range("C41").Value = 14
$ws.Range("D41").Value = 55
$ws.Range("E41").Value = 9
$ws.Range("A42").Value = "Badajoz"
$ws.Range("B42").Value = 111
$ws.Range("C42").Value = 5
$ws.Range("D42").Value = 104
$ws.Range("E42").Value = 2
$ws.Range("A43").Value = "Castello/Castellon"
$ws.Range("B43").Value = 104
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = 99
$ws.Range("E43").Value = 4
$ws.Range("A44").Value = "Ourense"
$ws.Range("B44").Value = 95
$ws.Range("C44").Value = 5
$ws.Range("D44").Value = 74
$ws.Range("E44").Value = 0
$ws.Range("A45").Value = "Zamora"
$ws.Range("B45").Value = 59
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 42
$ws.Range("E45").Value = 3
$ws.Range("A46").Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Range("B46").Value = 58
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 58
$ws.Range("E46").Value = 3
$ws.Range("A47").Value = "Almeria"
$ws.Range("B47").Value = 58
$ws.Range("C47").Value = 72
$ws.Range("D47").Value = 53
$ws.Range("E47").Value = 1
$ws.Range("A48").Value = "Lugo"
$ws.Range("B48").Value = 58
$ws.Range("C48").Value = 5
$ws.Range("D48").Value = 53
$ws.Range("E48").Value = 1
$ws.Range("A49").Value = "Teruel"
$ws.Range("B49").Value = 47
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 45
$ws.Range("E49").Value = 2
$ws.Range("A50").Value = "Huelva"
$ws.Range("B50").Value = 47
$ws.Range("C50").Value = 72
$ws.Range("D50").Value = 37
$ws.Range("E50").Value = 0
$ws.Range("A51").Value = "Palencia"
$ws.Range("B51").Value = 41
$ws.Range("C51").Value = 2
$ws.Range("D51").Value = 25
$ws.Range("E51").Value = 0
$ws.Range("A52").Value = "Huesca"
$ws.Range("B52").Value = 37
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 37
$ws.Range("E52").Value = 0
$ws.Range("A53").Value = "Melilla"
$ws.Range("B53").Value = 25
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 25
$ws.Range("E53").Value = 0
$ws.Range("A54").Value = "Ibiza"
$ws.Range("B54").Value = 21
$ws.Range("C54").Value = 10
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = 8
$ws.Range("A55").Value = "Menorca"
$ws.Range("B55").Value = 15
$ws.Range("C55").Value = 10
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 8
$ws.Range("A56").Value = "Fuerteventura"
$ws.Range("B56").Value = 14
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 13
$ws.Range("E56").Value = 0
$ws.Range("A57").Value = "La Palma"
$ws.Range("B57").Value = 14
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 11
$ws.Range("E57").Value = 0
$ws.Range("A58").Value = "Lanzarote"
$ws.Range("B58").Value = 9
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 4
$ws.Range("E58").Value = 0
$ws.Range("A59").Value = "Arroyo de la Luz"
$ws.Range("B59").Value = 7
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 7
$ws.Range("E59").Value = 0
$ws.Range("A60").Value = "Ceuta"
$ws.Range("B60").Value = 5
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 5
$ws.Range("E60").Value = 0
$ws.Range("A61").Value = "La Gomera"
$ws.Range("B61").Value = 3
$ws.Range("C61").Value = 2
$ws.Range("D61").Value = 1
$ws.Range("E61").Value = 0
$ws.Range("A62").Value = "El Hierro"
$ws.Range("B62").Value = 2
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = 0
$ws.Range("A63").Value = "Formentera"
$ws.Range("B63").Value = 0
$ws.Range("C63").Value = 10
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 8
